# Updated cryptos list on Tue Apr  2 09:29:47 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the crypto list on the active worksheet, and reorders the Bittensor/Maker
# rows (35 and 36) to reflect the coin's new rank.
#
# All of the Price/Volume cells are stored as plain text in the workbook
# (e.g. "66.300.33", "1.00", "  -4.59%  "), so before writing any value that
# Excel could otherwise auto-convert to a number (things like "1.00" or
# "560.49"), the target range is forced to Text format. The format/style is
# restored back to the sheet's default afterwards so no stray formatting is
# left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '66.300.33'
$ws.Range("E2").Value = '  -4.59%  '
$ws.Range("D3").Value = '3.341.47'
$ws.Range("E3").Value = '  -5.80%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '560.49'
$ws.Range("E5").Value = '  -4.23%  '
$ws.Range("D6").Value = '182.24'
$ws.Range("E6").Value = '  -7.83%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -3.02%  '
$ws.Range("D9").Value = '3.329.10'
$ws.Range("E9").Value = '  -5.77%  '
$ws.Range("E10").Value = '  -9.36%  '
$ws.Range("E11").Value = '  -6.75%  '
$ws.Range("D12").Value = '47.55'
$ws.Range("E12").Value = '  -8.45%  '
$ws.Range("D13").Value = '0.0000265'
$ws.Range("E13").Value = '  -7.68%  '
$ws.Range("D14").Value = '8.66'
$ws.Range("E14").Value = '  -6.52%  '
$ws.Range("D15").Value = '3.879.41'
$ws.Range("E15").Value = '  -5.56%  '
$ws.Range("D16").Value = '606.73'
$ws.Range("E16").Value = '  -8.60%  '
$ws.Range("D17").Value = '66.441.59'
$ws.Range("E17").Value = '  -4.59%  '
$ws.Range("D18").Value = '18.05'
$ws.Range("E18").Value = '  -2.78%  '
$ws.Range("D19").Value = '3.337.77'
$ws.Range("E19").Value = '  -5.85%  '
$ws.Range("E20").Value = '  -3.72%  '
$ws.Range("D21").Value = '11.46'
$ws.Range("E21").Value = '  -7.85%  '
$ws.Range("D22").Value = '0.907'
$ws.Range("E22").Value = '  -6.41%  '
$ws.Range("D23").Value = '16.87'
$ws.Range("E23").Value = '  -7.62%  '
$ws.Range("D24").Value = '5.06'
$ws.Range("E24").Value = '  -4.33%  '
$ws.Range("D25").Value = '99.98'
$ws.Range("E25").Value = '  -5.17%  '
$ws.Range("D26").Value = '4.07'
$ws.Range("E26").Value = '  -7.03%  '
$ws.Range("D27").Value = '6.01'
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("E28").Value = '  -7.90%  '
$ws.Range("D29").Value = '9.35'
$ws.Range("E29").Value = '  -8.18%  '
$ws.Range("E30").Value = '  -9.83%  '
$ws.Range("D31").Value = '30.38'
$ws.Range("E31").Value = '  -9.30%  '
$ws.Range("E32").Value = '  -8.49%  '
$ws.Range("E33").Value = '  -14.79%  '
$ws.Range("D34").Value = '11.04'
$ws.Range("E34").Value = '  -7.10%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '3.868.64'
$ws.Range("E35").Value = '  +2.12%  '
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").Value = '547.45'
$ws.Range("E36").Value = '  +8.52%  '
$ws.Range("E37").Value = '  -5.54%  '
$ws.Range("D38").Value = '57.48'
$ws.Range("E38").Value = '  -7.16%  '
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").Value = '3.41'
$ws.Range("E40").Value = '  -8.70%  '
$ws.Range("D41").Value = '0.0₃0716'
$ws.Range("E41").Value = '  -12.02%  '
$ws.Range("D42").Value = '2.67'
$ws.Range("E42").Value = '  -9.14%  '
$ws.Range("E43").Value = '  -6.91%  '
$ws.Range("D44").Value = '0.343'
$ws.Range("E44").Value = '  -8.16%  '
$ws.Range("D45").Value = '32.13'
$ws.Range("E45").Value = '  -7.37%  '
$ws.Range("E46").Value = '  +18.09%  '
$ws.Range("E47").Value = '  -8.78%  '
$ws.Range("D48").Value = '3.11'
$ws.Range("E48").Value = '  -8.37%  '
$ws.Range("E49").Value = '  -8.94%  '
$ws.Range("E50").Value = '  -4.74%  '
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  -0.11%  '

$priceRange.NumberFormat = "General"
$priceRange.Style = "Normal"
